# This edit removes the extra "karnata Captial City" quiz question that had
# been added in row 7 (A7:F7), reverting the sheet back to the shorter
# question list. Clearing those cells also drops the only references to the
# 5 shared strings it used ("karnata Captial City", "Mandya", "Banglore ",
# "mysore ", "DK"), so they fall out of the shared-string table on save and
# every later shared-string index shifts down automatically (e.g. row 5's
# answer cell, which keeps the same text but now points at a lower index).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the now-unwanted question row; keep the existing cell styling.
$ws.Range("A7:F7").ClearContents()

# Row 6 had an explicit wrapped-text height to match the long row 7 entry;
# auto-fit it back down to the sheet's normal row height now that row 7 is
# short/empty again.
$ws.Rows(6).AutoFit()

# Update the active selection to reflect the edit: the cleared row, and the
# next empty row below the question list.
$ws.Range("A7:F7").Select() | Out-Null
$ws.Range("B17").Select() | Out-Null
